{"js": "// Apply the dated-worksheet update: new date + each three-digit x one-digit\n// multiplication equation's operands/result, cell-by-cell, via exact text search.\nconst replacements = [\n  [\"2025-08-02 Saturday\", \"2025-08-03 Sunday\"],\n  [\"967\u00d79=8703\", \"621\u00d79=5589\"],\n  [\"426\u00d77=2982\", \"766\u00d75=3830\"],\n  [\"702\u00d76=4212\", \"756\u00d73=2268\"],\n  [\"190\u00d78=1520\", \"479\u00d77=3353\"],\n  [\"400\u00d73=1200\", \"873\u00d74=3492\"],\n  [\"614\u00d77=4298\", \"358\u00d77=2506\"],\n  [\"722\u00d72=1444\", \"613\u00d76=3678\"],\n  [\"465\u00d73=1395\", \"990\u00d76=5940\"],\n  [\"490\u00d77=3430\", \"861\u00d74=3444\"],\n  [\"403\u00d79=3627\", \"550\u00d79=4950\"],\n  [\"479\u00d74=1916\", \"333\u00d77=2331\"],\n  [\"429\u00d74=1716\", \"297\u00d78=2376\"],\n  [\"999\u00d78=7992\", \"105\u00d78=840\"],\n  [\"388\u00d76=2328\", \"558\u00d72=1116\"],\n  [\"834\u00d76=5004\", \"724\u00d74=2896\"],\n  [\"982\u00d72=1964\", \"361\u00d78=2888\"],\n  [\"841\u00d77=5887\", \"573\u00d73=1719\"],\n  [\"916\u00d75=4580\", \"506\u00d74=2024\"],\n  [\"352\u00d73=1056\", \"884\u00d77=6188\"],\n  [\"775\u00d78=6200\", \"781\u00d75=3905\"],\n  [\"645\u00d72=1290\", \"801\u00d78=6408\"],\n  [\"933\u00d76=5598\", \"822\u00d74=3288\"],\n  [\"364\u00d73=1092\", \"914\u00d76=5484\"],\n  [\"343\u00d77=2401\", \"559\u00d78=4472\"],\n  [\"124\u00d75=620\", \"783\u00d73=2349\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n\n", "ps1": "# Update the dated worksheet: refresh the header date and every\n# three-digit x one-digit multiplication equation cell in the table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2025-08-02 Saturday', '2025-08-03 Sunday'),\n    @('967\u00d79=8703', '621\u00d79=5589'),\n    @('426\u00d77=2982', '766\u00d75=3830'),\n    @('702\u00d76=4212', '756\u00d73=2268'),\n    @('190\u00d78=1520', '479\u00d77=3353'),\n    @('400\u00d73=1200', '873\u00d74=3492'),\n    @('614\u00d77=4298', '358\u00d77=2506'),\n    @('722\u00d72=1444', '613\u00d76=3678'),\n    @('465\u00d73=1395', '990\u00d76=5940'),\n    @('490\u00d77=3430', '861\u00d74=3444'),\n    @('403\u00d79=3627', '550\u00d79=4950'),\n    @('479\u00d74=1916', '333\u00d77=2331'),\n    @('429\u00d74=1716', '297\u00d78=2376'),\n    @('999\u00d78=7992', '105\u00d78=840'),\n    @('388\u00d76=2328', '558\u00d72=1116'),\n    @('834\u00d76=5004', '724\u00d74=2896'),\n    @('982\u00d72=1964', '361\u00d78=2888'),\n    @('841\u00d77=5887', '573\u00d73=1719'),\n    @('916\u00d75=4580', '506\u00d74=2024'),\n    @('352\u00d73=1056', '884\u00d77=6188'),\n    @('775\u00d78=6200', '781\u00d75=3905'),\n    @('645\u00d72=1290', '801\u00d78=6408'),\n    @('933\u00d76=5598', '822\u00d74=3288'),\n    @('364\u00d73=1092', '914\u00d76=5484'),\n    @('343\u00d77=2401', '559\u00d78=4472'),\n    @('124\u00d75=620', '783\u00d73=2349'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $ok = $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $ok) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n"}
